# Commit: "Removing TWB PNPC table and adding two extra fields to TWB
# Episode instead."
#
# 1) Delete the "TWB PNPCs" worksheet entirely.
# 2) Add two new fields to "TWB Episodes", inserted right after the
#    twb_primary_nominated_professional_consent_date column:
#      - twb_primary_nominated_professional_contact_entry_date
#      - twb_primary_nominated_professional_contact_exit_date

$wb = $excel.ActiveWorkbook
[void]($excel.DisplayAlerts = $false)

# --- 1. Remove the "TWB PNPCs" sheet -------------------------------------
$pnpc = $wb.Worksheets.Item("TWB PNPCs")
[void]$pnpc.Delete()

# --- 2. Insert the two new PNP contact date columns into "TWB Episodes" --
$twbEpisodes = $wb.Worksheets.Item("TWB Episodes")

# Insert two blank columns at L:M, pushing the existing
# twb_previous_suicide_attempts / twb_method_of_suicide_attempt columns
# (and their data) two slots to the right, into N:O.
[void]$twbEpisodes.Columns("L:M").Insert()

# New headers for the freshly inserted L:M columns.
$twbEpisodes.Range("L1").Value2 = "twb_primary_nominated_professional_contact_entry_date"
$twbEpisodes.Range("M1").Value2 = "twb_primary_nominated_professional_contact_exit_date"

# New data values for the two inserted columns.
$twbEpisodes.Range("L2").Value2 = 16042020
$twbEpisodes.Range("M2").Value2 = 9099999

$twbEpisodes.Range("L3").Value2 = 9099999
$twbEpisodes.Range("M3").Value2 = 9099999
